$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912" (sheet1) ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:26:12"
$ws1.Range("A3").Value = "Total filas: 2"

$ws1.Range("A6").Value = "02:26:12"
$ws1.Range("B6").Value = "02:58"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 32
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = "02:26:12"
$ws1.Range("B7").Value = "03:50"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 84
$ws1.Range("E7").Value = "LP1912"

# Row 8 no longer exists - remove it entirely
$ws1.Range("A8:E8").ClearContents()

# --- Sheet "LP1912-215" (sheet2) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:26:12"

$ws2.Range("A6").Value = "02:26:12"
$ws2.Range("D6").Value = 32

# --- Sheet "6203-6173" (sheet3) ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:26:12"
